$wb = $excel.ActiveWorkbook

function Set-PisoValues($ws) {
    $arr = New-Object 'object[,]' 30,20
    $arr[0,0] = 1
    $arr[0,1] = 11
    $arr[0,2] = 21
    $arr[0,3] = 31
    $arr[0,4] = 41
    $arr[0,5] = 51
    $arr[0,6] = 61
    $arr[0,7] = 71
    $arr[0,8] = 81
    $arr[0,9] = 91
    $arr[0,10] = 1
    $arr[0,11] = 11
    $arr[0,12] = 21
    $arr[0,13] = 31
    $arr[0,14] = 41
    $arr[0,15] = 51
    $arr[0,16] = 61
    $arr[0,17] = 71
    $arr[0,18] = 81
    $arr[0,19] = 91
    $arr[1,0] = 2
    $arr[1,9] = 92
    $arr[1,10] = 2
    $arr[1,19] = 92
    $arr[2,0] = 3
    $arr[2,9] = 93
    $arr[2,10] = 3
    $arr[2,19] = 93
    $arr[3,0] = 4
    $arr[3,9] = 94
    $arr[3,10] = 4
    $arr[3,19] = 94
    $arr[4,0] = 5
    $arr[4,9] = 95
    $arr[4,10] = 5
    $arr[4,19] = 95
    $arr[5,0] = 6
    $arr[5,9] = 96
    $arr[5,10] = 6
    $arr[5,19] = 96
    $arr[6,0] = 7
    $arr[6,9] = 97
    $arr[6,10] = 7
    $arr[6,19] = 97
    $arr[7,0] = 8
    $arr[7,9] = 98
    $arr[7,10] = 8
    $arr[7,19] = 98
    $arr[8,0] = 9
    $arr[8,9] = 99
    $arr[8,10] = 9
    $arr[8,19] = 99
    $arr[9,0] = 10
    $arr[9,1] = 20
    $arr[9,2] = 30
    $arr[9,3] = 40
    $arr[9,4] = 50
    $arr[9,5] = 60
    $arr[9,6] = 70
    $arr[9,7] = 80
    $arr[9,8] = 90
    $arr[9,9] = 100
    $arr[9,10] = 10
    $arr[9,11] = 20
    $arr[9,12] = 30
    $arr[9,13] = 40
    $arr[9,14] = 50
    $arr[9,15] = 60
    $arr[9,16] = 70
    $arr[9,17] = 80
    $arr[9,18] = 90
    $arr[9,19] = 100
    $arr[10,0] = 1
    $arr[10,1] = 11
    $arr[10,2] = 21
    $arr[10,3] = 31
    $arr[10,4] = 41
    $arr[10,5] = 51
    $arr[10,6] = 61
    $arr[10,7] = 71
    $arr[10,8] = 81
    $arr[10,9] = 91
    $arr[10,10] = 1
    $arr[10,11] = 11
    $arr[10,12] = 21
    $arr[10,13] = 31
    $arr[10,14] = 41
    $arr[10,15] = 51
    $arr[10,16] = 61
    $arr[10,17] = 71
    $arr[10,18] = 81
    $arr[10,19] = 91
    $arr[11,0] = 2
    $arr[11,9] = 92
    $arr[11,10] = 2
    $arr[11,19] = 92
    $arr[12,0] = 3
    $arr[12,9] = 93
    $arr[12,10] = 3
    $arr[12,19] = 93
    $arr[13,0] = 4
    $arr[13,9] = 94
    $arr[13,10] = 4
    $arr[13,19] = 94
    $arr[14,0] = 5
    $arr[14,9] = 95
    $arr[14,10] = 5
    $arr[14,19] = 95
    $arr[15,0] = 6
    $arr[15,9] = 96
    $arr[15,10] = 6
    $arr[15,19] = 96
    $arr[16,0] = 7
    $arr[16,9] = 97
    $arr[16,10] = 7
    $arr[16,19] = 97
    $arr[17,0] = 8
    $arr[17,9] = 98
    $arr[17,10] = 8
    $arr[17,19] = 98
    $arr[18,0] = 9
    $arr[18,9] = 99
    $arr[18,10] = 9
    $arr[18,19] = 99
    $arr[19,0] = 10
    $arr[19,1] = 20
    $arr[19,2] = 30
    $arr[19,3] = 40
    $arr[19,4] = 50
    $arr[19,5] = 60
    $arr[19,6] = 70
    $arr[19,7] = 80
    $arr[19,8] = 90
    $arr[19,9] = 100
    $arr[19,10] = 10
    $arr[19,11] = 20
    $arr[19,12] = 30
    $arr[19,13] = 40
    $arr[19,14] = 50
    $arr[19,15] = 60
    $arr[19,16] = 70
    $arr[19,17] = 80
    $arr[19,18] = 90
    $arr[19,19] = 100
    $arr[20,0] = 1
    $arr[20,1] = 11
    $arr[20,2] = 21
    $arr[20,3] = 31
    $arr[20,4] = 41
    $arr[20,5] = 51
    $arr[20,6] = 61
    $arr[20,7] = 71
    $arr[20,8] = 81
    $arr[20,9] = 91
    $arr[21,0] = 2
    $arr[21,9] = 92
    $arr[22,0] = 3
    $arr[22,9] = 93
    $arr[23,0] = 4
    $arr[23,9] = 94
    $arr[24,0] = 5
    $arr[24,9] = 95
    $arr[25,0] = 6
    $arr[25,9] = 96
    $arr[26,0] = 7
    $arr[26,9] = 97
    $arr[27,0] = 8
    $arr[27,9] = 98
    $arr[28,0] = 9
    $arr[28,9] = 99
    $arr[29,0] = 10
    $arr[29,1] = 20
    $arr[29,2] = 30
    $arr[29,3] = 40
    $arr[29,4] = 50
    $arr[29,5] = 60
    $arr[29,6] = 70
    $arr[29,7] = 80
    $arr[29,8] = 90
    $arr[29,9] = 100
    $ws.Range("A1:T30").Value = $arr
}

# Piso5 and Piso6 get the filled-in numeric grid (rows 1-20 cycle 1-10,
# columns A-J then repeat K-T; rows 21-30 only use columns A-J).
$ws5 = $wb.Worksheets.Item("Piso5")
Set-PisoValues $ws5

$ws6 = $wb.Worksheets.Item("Piso6")
Set-PisoValues $ws6

# View/selection updates.
# Piso5: zoom to 70%, scroll/select near J13 (was G29 at 100%).
$ws5.Activate()
$win = $wb.Windows.Item(1)
$win.Zoom = 70
$ws5.Range("J13").Select()

# Piso6 becomes the active/selected tab, zoom 70%, select K11:K20 (was F3 at 100%).
$ws6.Activate()
$win.Zoom = 70
$ws6.Range("K11:K20").Select()
